$d = $word.ActiveDocument

# --- Step 1: Split "Modera" into "Moderna" (with proofErr marks) and insert the new
# "Analysis of Covid-19 Deaths vs. Influenza Deaths" section right after that paragraph. ---
$p19 = $d.Paragraphs.Item(19)
$p20 = $d.Paragraphs.Item(20)
if ($p19.Range.Text -notlike "*Modera vaccines*") {
    throw "Unexpected document content at paragraph 19: $($p19.Range.Text)"
}
$combined = $d.Range($p19.Range.Start, $p20.Range.End)
$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">We also reviewed to see if there </w:t></w:r><w:r><w:t>were</w:t></w:r><w:r><w:t xml:space="preserve"> any relationships between the Pfizer and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Moder</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> vaccines. We found no clear relationships to the breakdown of distribution of the </w:t></w:r><w:r><w:t>manufacturers</w:t></w:r><w:r><w:t>. The only</w:t></w:r><w:r><w:t xml:space="preserve"> relationship was the total allocations to the total populations. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Analysis </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>f Covid-19 Deaths vs. Influenza Deaths</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">The data utilized was retrieved from the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">National Center for Health </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>and Statistics</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>NCHS</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Utilizing Line Graphs</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 1</w:t></w:r><w:r><w:t xml:space="preserve">: Evaluate </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Total Deaths</w:t></w:r><w:r><w:t xml:space="preserve"> for the past five years.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="0" w:name="_Hlk63291582"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Step </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve">: Evaluate </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Deaths</w:t></w:r><w:r><w:t xml:space="preserve"> for the past five years</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Step </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r><w:r><w:t xml:space="preserve">: Evaluate </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Covid-19</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Deaths</w:t></w:r><w:r><w:t xml:space="preserve"> for </w:t></w:r><w:r><w:t>2020</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="1" w:name="_Hlk63291761"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4a</w:t></w:r><w:r><w:t xml:space="preserve">: Calculate </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Covid-19 Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>in targeted 2020 timeframe.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="2" w:name="_Hlk63291835"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>b</w:t></w:r><w:r><w:t xml:space="preserve">: Calculate </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Covid-19 </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Non-</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">in </w:t></w:r><w:r><w:t xml:space="preserve">targeted </w:t></w:r><w:r><w:t>2020</w:t></w:r><w:r><w:t xml:space="preserve"> timeframe</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="3" w:name="_Hlk63291945"/><w:bookmarkEnd w:id="2"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>c</w:t></w:r><w:r><w:t xml:space="preserve">: Calculate </w:t></w:r><w:bookmarkStart w:id="4" w:name="_Hlk63292056"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Influenza</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="4"/><w:r><w:t xml:space="preserve">in </w:t></w:r><w:r><w:t xml:space="preserve">targeted </w:t></w:r><w:r><w:t>2020</w:t></w:r><w:r><w:t xml:space="preserve"> timeframe</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="5" w:name="_Hlk63292008"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>d</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>Evaluate</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Targeted </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">in </w:t></w:r><w:r><w:t xml:space="preserve">targeted </w:t></w:r><w:r><w:t>2020</w:t></w:r><w:r><w:t xml:space="preserve"> timeframe</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="6" w:name="_Hlk63292136"/><w:bookmarkEnd w:id="5"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">Use correlation analysis to verify </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Influenza</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>were statistically consistent with the same timeframe in previous years.</w:t></w:r><w:bookmarkEnd w:id="6"/></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>f</w:t></w:r><w:r><w:t xml:space="preserve">: Use </w:t></w:r><w:r><w:t>ANOVA</w:t></w:r><w:r><w:t xml:space="preserve"> analysis to verify </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Influenza</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia Deaths</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">were statistically </w:t></w:r><w:r><w:t>distributed</w:t></w:r><w:r><w:t xml:space="preserve"> with the same timeframe in previous years.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Utilizing Pie Charts</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:bookmarkStart w:id="7" w:name="_Hlk63292442"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 5</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Compare 2020 </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Influenza Pneumonia</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Non-Pneumonia Deaths</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkEnd w:id="7"/></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Step </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>6</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> Compare 2020 </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Covid-19 </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Pneumonia</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Non-Pneumonia Deaths</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>References</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 7</w:t></w:r><w:r><w:t xml:space="preserve">: Provide a link concerning </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Comorbidities</w:t></w:r><w:r><w:t xml:space="preserve"> associated with Covid-19 deaths.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Step </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>8</w:t></w:r><w:r><w:t xml:space="preserve">: Provide a link concerning </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>C</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">ARES Act </w:t></w:r><w:r><w:t>compensation</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>made available to hospitals during the Covid-19 epidemic</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Thoughts</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 9</w:t></w:r><w:r><w:t>: Reflect on a little of what was seen in the previous charts.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Step 10</w:t></w:r><w:r><w:t>: Provide some explanation of link information.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Step 11</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Perform calculations to obtain an estimated number of </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Weighted Covid-19 Deaths</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/></w:rPr><w:t>A Final Question</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Step 12: Pose a final question concerning the results, calculations, and information reviewed</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="24"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>
'@
$combined.InsertXML($fragment)

Write-Host "done"
